function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $ws.Range("H2").Style
}

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)   # 总计
$template = $wb.Worksheets.Item(2) # currently 2022-Q3, used as a style/layout template

# Insert a new sheet (copy of the Q3 template) right after 总计; this will hold 2022-Q4 data
$template.Copy($null, $sheet1)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q4"

# Populate 2022-Q4 sheet with its own data (preserve exact text formatting via Set-TextValue)
Set-TextValue $ws "B2" "710001"
Set-TextValue $ws "C2" "富安达优势成长混合"
Set-TextValue $ws "D2" "8.31"
Set-TextValue $ws "E2" "93.43"
Set-TextValue $ws "F2" "4.76"
Set-TextValue $ws "G2" "0.3956"
$ws.Range("H2").Value = 5
Set-TextValue $ws "B3" "163804"
Set-TextValue $ws "C3" "中银收益混合A"
Set-TextValue $ws "D3" "15.84"
Set-TextValue $ws "E3" "88.15"
Set-TextValue $ws "F3" "2.41"
Set-TextValue $ws "G3" "0.3817"
$ws.Range("H3").Value = 9
Set-TextValue $ws "B4" "163822"
Set-TextValue $ws "C4" "中银主题策略混合A"
Set-TextValue $ws "D4" "17.43"
Set-TextValue $ws "E4" "87.27"
Set-TextValue $ws "F4" "2.11"
Set-TextValue $ws "G4" "0.3678"
$ws.Range("H4").Value = 10
Set-TextValue $ws "B5" "014505"
Set-TextValue $ws "C5" "中银收益混合C"
Set-TextValue $ws "D5" "6.72"
Set-TextValue $ws "E5" "88.15"
Set-TextValue $ws "F5" "2.41"
Set-TextValue $ws "G5" "0.1620"
$ws.Range("H5").Value = 9
Set-TextValue $ws "B6" "014103"
Set-TextValue $ws "C6" "富安达成长价值一年持有期混合A"
Set-TextValue $ws "D6" "2.30"
Set-TextValue $ws "E6" "89.75"
Set-TextValue $ws "F6" "5.19"
Set-TextValue $ws "G6" "0.1194"
$ws.Range("H6").Value = 4
Set-TextValue $ws "B7" "015386"
Set-TextValue $ws "C7" "中银主题策略混合C"
Set-TextValue $ws "D7" "4.86"
Set-TextValue $ws "E7" "87.27"
Set-TextValue $ws "F7" "2.11"
Set-TextValue $ws "G7" "0.1025"
$ws.Range("H7").Value = 10
Set-TextValue $ws "B8" "010965"
Set-TextValue $ws "C8" "中银鑫新消费成长混合A"
Set-TextValue $ws "D8" "3.24"
Set-TextValue $ws "E8" "88.64"
Set-TextValue $ws "F8" "2.22"
Set-TextValue $ws "G8" "0.0719"
$ws.Range("H8").Value = 10
Set-TextValue $ws "B9" "009789"
Set-TextValue $ws "C9" "富安达科技创新混合"
Set-TextValue $ws "D9" "0.43"
Set-TextValue $ws "E9" "93.91"
Set-TextValue $ws "F9" "4.52"
Set-TextValue $ws "G9" "0.0194"
$ws.Range("H9").Value = 9
Set-TextValue $ws "B10" "010962"
Set-TextValue $ws "C10" "中银鑫新消费成长混合C"
Set-TextValue $ws "D10" "0.69"
Set-TextValue $ws "E10" "88.64"
Set-TextValue $ws "F10" "2.22"
Set-TextValue $ws "G10" "0.0153"
$ws.Range("H10").Value = 10
Set-TextValue $ws "B11" "014104"
Set-TextValue $ws "C11" "富安达成长价值一年持有期混合C"
Set-TextValue $ws "D11" "0.22"
Set-TextValue $ws "E11" "89.75"
Set-TextValue $ws "F11" "5.19"
Set-TextValue $ws "G11" "0.0114"
$ws.Range("H11").Value = 4
Set-TextValue $ws "B12" "960012"
Set-TextValue $ws "C12" "中银收益混合H"
Set-TextValue $ws "D12" "0.04"
Set-TextValue $ws "E12" "88.15"
Set-TextValue $ws "F12" "2.41"
Set-TextValue $ws "G12" "0.0010"
$ws.Range("H12").Value = 9

# Restore the active/selected tab to the last sheet (2020-Q4), matching original workbook state
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

# Update the 总计 (summary) sheet: insert a new row for 2022-Q4 and shift existing rows down
$sheet1.Rows.Item(2).Insert()
$sheet1.Range("B2:D2").Style = $sheet1.Range("B3").Style

$sheet1.Range("A2").Value = 0
$sheet1.Range("B2").Value = "2022-Q4"
$sheet1.Range("C2").Value = 11
$sheet1.Range("D2").Value = 1.65
$sheet1.Range("A3").Value = 1
$sheet1.Range("B3").Value = "2022-Q3"
$sheet1.Range("C3").Value = 11
$sheet1.Range("D3").Value = 2.13
$sheet1.Range("A4").Value = 2
$sheet1.Range("B4").Value = "2022-Q2"
$sheet1.Range("C4").Value = 3
$sheet1.Range("D4").Value = 0.44
$sheet1.Range("A5").Value = 3
$sheet1.Range("B5").Value = "2022-Q1"
$sheet1.Range("C5").Value = 2
$sheet1.Range("D5").Value = 0.38
$sheet1.Range("A6").Value = 4
$sheet1.Range("B6").Value = "2021-Q4"
$sheet1.Range("C6").Value = 8
$sheet1.Range("D6").Value = 2.43
$sheet1.Range("A7").Value = 5
$sheet1.Range("B7").Value = "2021-Q3"
$sheet1.Range("C7").Value = 3
$sheet1.Range("D7").Value = 0.51
$sheet1.Range("A8").Value = 6
$sheet1.Range("B8").Value = "2021-Q2"
$sheet1.Range("C8").Value = 1
$sheet1.Range("D8").Value = 0.4
$sheet1.Range("A9").Value = 7
$sheet1.Range("B9").Value = "2021-Q1"
$sheet1.Range("C9").Value = 6
$sheet1.Range("D9").Value = 1
$sheet1.Range("A10").Value = 8
$sheet1.Range("B10").Value = "2020-Q4"
$sheet1.Range("C10").Value = 1
$sheet1.Range("D10").Value = 0.61
$sheet1.Range("A2").Style = $sheet1.Range("A3").Style

Write-Host "Done"